# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.284.50'
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").Value = '1.650.34'
$ws.Range("E3").Value = '  -0.53%  '

$ws.Range("E4").Value = '  -0.15%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '218.63'
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("E6").Value = '  +1.70%  '

$ws.Range("E8").Value = '  +0.85%  '

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0628'
$ws.Range("E9").Value = '  +0.18%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '20.22'
$ws.Range("E10").Value = '  +2.14%  '

$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").Value = '1.879.95'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("D13").Value = '1.637.20'
$ws.Range("E13").Value = '  -1.26%  '

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '4.14'
$ws.Range("E14").Value = '  -1.45%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '0.539'
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("E16").Value = '  +1.90%  '

$ws.Range("D17").Value = '27.242.79'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").Value = '0.0₃0743'
$ws.Range("E18").Value = '  +0.61%  '

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '221.18'
$ws.Range("E19").Value = '  -0.16%  '

$ws.Range("E20").Value = '  -0.18%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '6.77'
$ws.Range("E21").Value = '  +0.73%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '4.46'
$ws.Range("E22").Value = '  +0.32%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '2.51'
$ws.Range("E23").Value = '  +2.98%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '9.26'
$ws.Range("E24").Value = '  -0.38%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '148.60'
$ws.Range("E25").Value = '  +0.80%  '

$ws.Range("E26").Value = '  -0.05%  '

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '7.42'
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("E28").Value = '  +0.21%  '

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '15.84'
$ws.Range("E29").Value = '  -0.73%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0507'
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("E31").Value = '  -0.44%  '

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '3.37'
$ws.Range("E32").Value = '  -0.08%  '

$ws.Range("E33").Value = '  +1.23%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.59'
$ws.Range("E34").Value = '  +1.28%  '

$ws.Range("D35").Value = '1.277.63'
$ws.Range("E35").Value = '  +0.42%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '2.46'
$ws.Range("E36").Value = '  +0.30%  '

$ws.Range("E37").Value = '  +1.22%  '

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.544'
$ws.Range("E38").Value = '  +0.92%  '

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.849'
$ws.Range("E39").Value = '  +2.18%  '

$ws.Range("E40").Value = '  -0.14%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.813'
$ws.Range("E41").Value = '  +0.51%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '2.22'
$ws.Range("E42").Value = '  +4.99%  '

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '5.38'
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").Value = '1.790.33'
$ws.Range("E44").Value = '  -0.59%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '62.83'
$ws.Range("E45").Value = '  +1.51%  '

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '92.57'
$ws.Range("E46").Value = '  -0.19%  '

$ws.Range("E47").Value = '  -0.89%  '

$ws.Range("D48").Value = '0.0₆0108'
$ws.Range("E48").Value = '  +19.12%  '

$ws.Range("E49").Value = '  -0.65%  '

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '7.71'
$ws.Range("E50").Value = '  +0.59%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0975'
$ws.Range("E51").Value = '  -0.43%  '
